{"js": "const replacements = [\n  [\"376\u00d73=1128\", \"416\u00d77=2912\"],\n  [\"429\u00d72=858\", \"937\u00d79=8433\"],\n  [\"840\u00d77=5880\", \"835\u00d74=3340\"],\n  [\"907\u00d74=3628\", \"450\u00d75=2250\"],\n  [\"246\u00d76=1476\", \"245\u00d72=490\"],\n  [\"603\u00d75=3015\", \"178\u00d73=534\"],\n  [\"466\u00d78=3728\", \"443\u00d77=3101\"],\n  [\"559\u00d75=2795\", \"234\u00d74=936\"],\n  [\"666\u00d72=1332\", \"894\u00d73=2682\"],\n  [\"132\u00d73=396\", \"602\u00d75=3010\"],\n  [\"978\u00d76=5868\", \"807\u00d78=6456\"],\n  [\"845\u00d77=5915\", \"636\u00d74=2544\"],\n  [\"214\u00d75=1070\", \"606\u00d79=5454\"],\n  [\"271\u00d78=2168\", \"529\u00d74=2116\"],\n  [\"855\u00d77=5985\", \"705\u00d72=1410\"],\n  [\"989\u00d76=5934\", \"426\u00d75=2130\"],\n  [\"461\u00d77=3227\", \"903\u00d75=4515\"],\n  [\"474\u00d77=3318\", \"311\u00d79=2799\"],\n  [\"959\u00d78=7672\", \"200\u00d74=800\"],\n  [\"389\u00d79=3501\", \"828\u00d72=1656\"],\n  [\"822\u00d78=6576\", \"172\u00d79=1548\"],\n  [\"530\u00d73=1590\", \"122\u00d78=976\"],\n  [\"840\u00d76=5040\", \"846\u00d78=6768\"],\n  [\"223\u00d73=669\", \"560\u00d73=1680\"],\n  [\"941\u00d77=6587\", \"703\u00d78=5624\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${oldText}`);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"376\u00d73=1128\"; New=\"416\u00d77=2912\"}\n    @{Old=\"429\u00d72=858\"; New=\"937\u00d79=8433\"}\n    @{Old=\"840\u00d77=5880\"; New=\"835\u00d74=3340\"}\n    @{Old=\"907\u00d74=3628\"; New=\"450\u00d75=2250\"}\n    @{Old=\"246\u00d76=1476\"; New=\"245\u00d72=490\"}\n    @{Old=\"603\u00d75=3015\"; New=\"178\u00d73=534\"}\n    @{Old=\"466\u00d78=3728\"; New=\"443\u00d77=3101\"}\n    @{Old=\"559\u00d75=2795\"; New=\"234\u00d74=936\"}\n    @{Old=\"666\u00d72=1332\"; New=\"894\u00d73=2682\"}\n    @{Old=\"132\u00d73=396\"; New=\"602\u00d75=3010\"}\n    @{Old=\"978\u00d76=5868\"; New=\"807\u00d78=6456\"}\n    @{Old=\"845\u00d77=5915\"; New=\"636\u00d74=2544\"}\n    @{Old=\"214\u00d75=1070\"; New=\"606\u00d79=5454\"}\n    @{Old=\"271\u00d78=2168\"; New=\"529\u00d74=2116\"}\n    @{Old=\"855\u00d77=5985\"; New=\"705\u00d72=1410\"}\n    @{Old=\"989\u00d76=5934\"; New=\"426\u00d75=2130\"}\n    @{Old=\"461\u00d77=3227\"; New=\"903\u00d75=4515\"}\n    @{Old=\"474\u00d77=3318\"; New=\"311\u00d79=2799\"}\n    @{Old=\"959\u00d78=7672\"; New=\"200\u00d74=800\"}\n    @{Old=\"389\u00d79=3501\"; New=\"828\u00d72=1656\"}\n    @{Old=\"822\u00d78=6576\"; New=\"172\u00d79=1548\"}\n    @{Old=\"530\u00d73=1590\"; New=\"122\u00d78=976\"}\n    @{Old=\"840\u00d76=5040\"; New=\"846\u00d78=6768\"}\n    @{Old=\"223\u00d73=669\"; New=\"560\u00d73=1680\"}\n    @{Old=\"941\u00d77=6587\"; New=\"703\u00d78=5624\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
